# Auto-generated edit script: updates Leve market-price snapshot values
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns H:N)
# across the ALC, ARM, CUL, GSM, LTW, WVR leve-profit sheets to match the
# refreshed Universalis market data snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 30000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
$ws.Range("H64").Value = 3948.5715
$ws.Range("I64").Value = 3815.2942
$ws.Range("J64").Value = 4154.5454
$ws.Range("K64").Value = 3815.2942
$ws.Range("L64").Value = 4154.5454
$ws.Range("M64").Value = -3567.2942
$ws.Range("N64").Value = -4650.5454
$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 30000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
$ws.Range("H67").Value = 3948.5715
$ws.Range("I67").Value = 3815.2942
$ws.Range("J67").Value = 4154.5454
$ws.Range("K67").Value = 3815.2942
$ws.Range("L67").Value = 4154.5454
$ws.Range("M67").Value = -2957.2942
$ws.Range("N67").Value = -5870.5454
$ws.Range("H69").Value = 12741.5
$ws.Range("I69").Value = 1675.3334
$ws.Range("J69").Value = 17484.143
$ws.Range("K69").Value = 5026.0002
$ws.Range("L69").Value = 52452.429
$ws.Range("M69").Value = -4152.0002
$ws.Range("N69").Value = -54200.429
$ws.Range("H70").Value = 1662.7142
$ws.Range("I70").Value = 1459.75
$ws.Range("J70").Value = 1743.9
$ws.Range("K70").Value = 4379.25
$ws.Range("L70").Value = 5231.700000000001
$ws.Range("M70").Value = -4109.25
$ws.Range("N70").Value = -5771.700000000001
$ws.Range("H72").Value = 12741.5
$ws.Range("I72").Value = 1675.3334
$ws.Range("J72").Value = 17484.143
$ws.Range("K72").Value = 15078.0006
$ws.Range("L72").Value = 157357.287
$ws.Range("M72").Value = -10710.0006
$ws.Range("N72").Value = -166093.287
$ws.Range("H73").Value = 1662.7142
$ws.Range("I73").Value = 1459.75
$ws.Range("J73").Value = 1743.9
$ws.Range("K73").Value = 4379.25
$ws.Range("L73").Value = 5231.700000000001
$ws.Range("M73").Value = -3443.25
$ws.Range("N73").Value = -7103.700000000001
$ws.Range("H74").Value = 4942.2856
$ws.Range("I74").Value = 4932.6665
$ws.Range("J74").Value = 4949.5
$ws.Range("K74").Value = 4932.6665
$ws.Range("L74").Value = 4949.5
$ws.Range("M74").Value = -3996.6665
$ws.Range("N74").Value = -6821.5
$ws.Range("H76").Value = 3057.7778
$ws.Range("I76").Value = 3035.8975
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3035.8975
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -2720.8975
$ws.Range("H77").Value = 4942.2856
$ws.Range("I77").Value = 4932.6665
$ws.Range("J77").Value = 4949.5
$ws.Range("K77").Value = 24663.3325
$ws.Range("L77").Value = 24747.5
$ws.Range("M77").Value = -19983.3325
$ws.Range("N77").Value = -34107.5
$ws.Range("H79").Value = 3057.7778
$ws.Range("I79").Value = 3035.8975
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3035.8975
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -1943.8975
$ws.Range("H80").Value = 3389.2222
$ws.Range("I80").Value = 2224.15
$ws.Range("J80").Value = 4074.5588
$ws.Range("K80").Value = 6672.450000000001
$ws.Range("L80").Value = 12223.6764
$ws.Range("M80").Value = -5674.450000000001
$ws.Range("N80").Value = -14219.6764
$ws.Range("H82").Value = 796.4167
$ws.Range("I82").Value = 796.4167
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2389.2501
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1983.2501
$ws.Range("H83").Value = 3389.2222
$ws.Range("I83").Value = 2224.15
$ws.Range("J83").Value = 4074.5588
$ws.Range("K83").Value = 20017.35
$ws.Range("L83").Value = 36671.0292
$ws.Range("M83").Value = -15025.35
$ws.Range("N83").Value = -46655.0292
$ws.Range("H85").Value = 796.4167
$ws.Range("I85").Value = 796.4167
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2389.2501
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -985.2501000000002
$ws.Range("H87").Value = 28999.75
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 28999.75
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 28999.75
$ws.Range("N87").Value = -31495.75
$ws.Range("H88").Value = 2388.7812
$ws.Range("I88").Value = 1683.1
$ws.Range("J88").Value = 2709.5454
$ws.Range("K88").Value = 1683.1
$ws.Range("L88").Value = 2709.5454
$ws.Range("M88").Value = -1277.1
$ws.Range("N88").Value = -3521.5454
$ws.Range("H90").Value = 28999.75
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 28999.75
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 86999.25
$ws.Range("N90").Value = -99479.25
$ws.Range("H91").Value = 2388.7812
$ws.Range("I91").Value = 1683.1
$ws.Range("J91").Value = 2709.5454
$ws.Range("K91").Value = 1683.1
$ws.Range("L91").Value = 2709.5454
$ws.Range("M91").Value = -279.0999999999999
$ws.Range("N91").Value = -5517.5454

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3410.7258
$ws.Range("I32").Value = 3055.4807
$ws.Range("J32").Value = 5258
$ws.Range("K32").Value = 3055.4807
$ws.Range("L32").Value = 5258
$ws.Range("M32").Value = -2768.4807
$ws.Range("N32").Value = -5832
$ws.Range("H61").Value = 1460.2858
$ws.Range("I61").Value = 1345.3334
$ws.Range("J61").Value = 2150
$ws.Range("K61").Value = 1345.3334
$ws.Range("L61").Value = 2150
$ws.Range("M61").Value = -1133.3334
$ws.Range("N61").Value = -2574
$ws.Range("H74").Value = 123270.78
$ws.Range("I74").Value = 167820.5
$ws.Range("J74").Value = 34171.332
$ws.Range("K74").Value = 167820.5
$ws.Range("L74").Value = 34171.332
$ws.Range("M74").Value = -166946.5
$ws.Range("N74").Value = -35919.332
$ws.Range("H77").Value = 123270.78
$ws.Range("I77").Value = 167820.5
$ws.Range("J77").Value = 34171.332
$ws.Range("K77").Value = 839102.5
$ws.Range("L77").Value = 170856.66
$ws.Range("M77").Value = -834734.5
$ws.Range("N77").Value = -179592.66
$ws.Range("H136").Value = 1460.2858
$ws.Range("I136").Value = 1345.3334
$ws.Range("J136").Value = 2150
$ws.Range("K136").Value = 4036.0002
$ws.Range("L136").Value = 6450
$ws.Range("M136").Value = -1486.0002
$ws.Range("N136").Value = -11550

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 196.28572
$ws.Range("I33").Value = 363.33334
$ws.Range("J33").Value = 71
$ws.Range("K33").Value = 2180.00004
$ws.Range("L33").Value = 426
$ws.Range("M33").Value = -1897.00004
$ws.Range("N33").Value = -992
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H63").Value = 5933.3335
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 5933.3335
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 17800.0005
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -19298.0005
$ws.Range("H64").Value = 6828.5
$ws.Range("I64").Value = 6500
$ws.Range("J64").Value = 6938
$ws.Range("K64").Value = 19500
$ws.Range("L64").Value = 20814
$ws.Range("M64").Value = -19230
$ws.Range("N64").Value = -21354
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H66").Value = 5933.3335
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 5933.3335
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 53400.0015
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -60888.0015
$ws.Range("H67").Value = 6828.5
$ws.Range("I67").Value = 6500
$ws.Range("J67").Value = 6938
$ws.Range("K67").Value = 19500
$ws.Range("L67").Value = 20814
$ws.Range("M67").Value = -18564
$ws.Range("N67").Value = -22686
$ws.Range("H68").Value = 1500
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4500
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6122
$ws.Range("H69").Value = 3202.2
$ws.Range("I69").Value = 2003.6666
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 6010.9998
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -5199.9998
$ws.Range("H70").Value = 6007.9
$ws.Range("I70").Value = 4869.75
$ws.Range("J70").Value = 6766.6665
$ws.Range("K70").Value = 14609.25
$ws.Range("L70").Value = 20299.9995
$ws.Range("M70").Value = -14294.25
$ws.Range("N70").Value = -20929.9995
$ws.Range("H71").Value = 1500
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 13500
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -21612
$ws.Range("H72").Value = 3202.2
$ws.Range("I72").Value = 2003.6666
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 18032.9994
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -13976.9994
$ws.Range("H73").Value = 6007.9
$ws.Range("I73").Value = 4869.75
$ws.Range("J73").Value = 6766.6665
$ws.Range("K73").Value = 14609.25
$ws.Range("L73").Value = 20299.9995
$ws.Range("M73").Value = -13517.25
$ws.Range("N73").Value = -22483.9995
$ws.Range("H74").Value = 8289.666999999999
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8289.666999999999
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 24869.001
$ws.Range("N74").Value = -26991.001
$ws.Range("H75").Value = 4673.7144
$ws.Range("I75").Value = 3050
$ws.Range("J75").Value = 4944.3335
$ws.Range("K75").Value = 9150
$ws.Range("L75").Value = 14833.0005
$ws.Range("M75").Value = -8152
$ws.Range("N75").Value = -16829.0005
$ws.Range("H76").Value = 4555
$ws.Range("I76").Value = 1515
$ws.Range("J76").Value = 5163
$ws.Range("K76").Value = 4545
$ws.Range("L76").Value = 15489
$ws.Range("M76").Value = -4162
$ws.Range("N76").Value = -16255
$ws.Range("H77").Value = 8289.666999999999
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8289.666999999999
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 74607.003
$ws.Range("N77").Value = -85215.003
$ws.Range("H78").Value = 4673.7144
$ws.Range("I78").Value = 3050
$ws.Range("J78").Value = 4944.3335
$ws.Range("K78").Value = 27450
$ws.Range("L78").Value = 44499.0015
$ws.Range("M78").Value = -22458
$ws.Range("N78").Value = -54483.0015
$ws.Range("H79").Value = 4555
$ws.Range("I79").Value = 1515
$ws.Range("J79").Value = 5163
$ws.Range("K79").Value = 4545
$ws.Range("L79").Value = 15489
$ws.Range("M79").Value = -3219
$ws.Range("N79").Value = -18141
$ws.Range("H80").Value = 2454.2856
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2454.2856
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 7362.8568
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -9234.856800000001
$ws.Range("H81").Value = 3407.8572
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 3809.1667
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 11427.5001
$ws.Range("M81").Value = -1877
$ws.Range("N81").Value = -13673.5001
$ws.Range("H82").Value = 2753.25
$ws.Range("I82").Value = 1013
$ws.Range("J82").Value = 3333.3333
$ws.Range("K82").Value = 3039
$ws.Range("L82").Value = 9999.999899999999
$ws.Range("M82").Value = -2633
$ws.Range("N82").Value = -10811.9999
$ws.Range("H83").Value = 2454.2856
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2454.2856
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 22088.5704
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -31448.5704
$ws.Range("H84").Value = 3407.8572
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 3809.1667
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 34282.5003
$ws.Range("M84").Value = -3384
$ws.Range("N84").Value = -45514.5003
$ws.Range("H85").Value = 2753.25
$ws.Range("I85").Value = 1013
$ws.Range("J85").Value = 3333.3333
$ws.Range("K85").Value = 3039
$ws.Range("L85").Value = 9999.999899999999
$ws.Range("M85").Value = -1635
$ws.Range("N85").Value = -12807.9999
$ws.Range("H86").Value = 384.2
$ws.Range("I86").Value = 380.66666
$ws.Range("J86").Value = 389.5
$ws.Range("K86").Value = 1141.99998
$ws.Range("L86").Value = 1168.5
$ws.Range("M86").Value = 44.00001999999995
$ws.Range("N86").Value = -3540.5
$ws.Range("H87").Value = 5000
$ws.Range("I87").Value = 5000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 15000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -13752
$ws.Range("H88").Value = 5600
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5600
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 16800
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -17656
$ws.Range("H89").Value = 384.2
$ws.Range("I89").Value = 380.66666
$ws.Range("J89").Value = 389.5
$ws.Range("K89").Value = 3425.99994
$ws.Range("L89").Value = 3505.5
$ws.Range("M89").Value = 2502.00006
$ws.Range("N89").Value = -15361.5
$ws.Range("H90").Value = 5000
$ws.Range("I90").Value = 5000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 45000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -38760
$ws.Range("H91").Value = 5600
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5600
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 16800
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -19764
$ws.Range("H131").Value = 895.3
$ws.Range("I131").Value = 600
$ws.Range("J131").Value = 898.28284
$ws.Range("K131").Value = 1800
$ws.Range("L131").Value = 2694.84852
$ws.Range("M131").Value = 3240
$ws.Range("N131").Value = -12774.84852

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 14000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 14000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 14000
$ws.Range("N15").Value = -14576
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 30000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 30000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 30000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 30000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H68").Value = 30000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 30000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31622
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H71").Value = 30000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 30000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -98112
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 2695.4092
$ws.Range("I80").Value = 2594.7368
$ws.Range("J80").Value = 3333
$ws.Range("K80").Value = 2594.7368
$ws.Range("L80").Value = 3333
$ws.Range("M80").Value = -1596.7368
$ws.Range("H81").Value = 14000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 14000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 14000
$ws.Range("N81").Value = -15996
$ws.Range("H82").Value = 30000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 30000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766
$ws.Range("H83").Value = 2695.4092
$ws.Range("I83").Value = 2594.7368
$ws.Range("J83").Value = 3333
$ws.Range("K83").Value = 12973.684
$ws.Range("L83").Value = 16665
$ws.Range("M83").Value = -7981.684000000001
$ws.Range("H84").Value = 14000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 14000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 42000
$ws.Range("N84").Value = -51984
$ws.Range("H85").Value = 30000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 30000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("H132").Value = 2409.818
$ws.Range("I132").Value = 2217.3333
$ws.Range("J132").Value = 2640.8
$ws.Range("K132").Value = 6651.999899999999
$ws.Range("L132").Value = 7922.400000000001
$ws.Range("M132").Value = -4121.999899999999
$ws.Range("N132").Value = -12982.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1540.2
$ws.Range("I16").Value = 1540.2
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1540.2
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1370.2
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H64").Value = 31333.334
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 31333.334
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 31333.334
$ws.Range("N64").Value = -31783.334
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H67").Value = 31333.334
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 31333.334
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 31333.334
$ws.Range("N67").Value = -32893.334
$ws.Range("H68").Value = 14444.444
$ws.Range("I68").Value = 22500
$ws.Range("J68").Value = 4375
$ws.Range("K68").Value = 22500
$ws.Range("L68").Value = 4375
$ws.Range("M68").Value = -21751
$ws.Range("N68").Value = -5873
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H71").Value = 14444.444
$ws.Range("I71").Value = 22500
$ws.Range("J71").Value = 4375
$ws.Range("K71").Value = 112500
$ws.Range("L71").Value = 21875
$ws.Range("M71").Value = -108756
$ws.Range("N71").Value = -29363
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H76").Value = 30000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 30000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30676
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H79").Value = 30000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 30000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32340
$ws.Range("H80").Value = 25000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 25000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -27246
$ws.Range("H81").Value = 30000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 30000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H82").Value = 3666.4443
$ws.Range("I82").Value = 1999.5
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 1999.5
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -1638.5
$ws.Range("N82").Value = -5722
$ws.Range("H83").Value = 25000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 25000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -86232
$ws.Range("H84").Value = 30000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 30000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H85").Value = 3666.4443
$ws.Range("I85").Value = 1999.5
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 1999.5
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -751.5
$ws.Range("N85").Value = -7496
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2742.7144
$ws.Range("I62").Value = 2600
$ws.Range("J62").Value = 2799.8
$ws.Range("K62").Value = 2600
$ws.Range("L62").Value = 2799.8
$ws.Range("M62").Value = -1976
$ws.Range("N62").Value = -4047.8
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 28000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 28000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 28000
$ws.Range("N64").Value = -28496
$ws.Range("H65").Value = 2742.7144
$ws.Range("I65").Value = 2600
$ws.Range("J65").Value = 2799.8
$ws.Range("K65").Value = 13000
$ws.Range("L65").Value = 13999
$ws.Range("M65").Value = -9880
$ws.Range("N65").Value = -20239
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 28000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 28000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 28000
$ws.Range("N67").Value = -29716
$ws.Range("H70").Value = 25932.777
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 25932.777
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 25932.777
$ws.Range("N70").Value = -26562.777
$ws.Range("H73").Value = 25932.777
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 25932.777
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 25932.777
$ws.Range("N73").Value = -28116.777
$ws.Range("H75").Value = 29933.334
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 29933.334
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 29933.334
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -31805.334
$ws.Range("H76").Value = 30000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 30000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630
$ws.Range("H78").Value = 29933.334
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 29933.334
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 89800.00199999999
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -99160.00199999999
$ws.Range("H79").Value = 30000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 30000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184
$ws.Range("H81").Value = 1087.75
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1087.75
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 2175.5
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -4297.5
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H84").Value = 1087.75
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1087.75
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 10877.5
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -21485.5
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H88").Value = 29175
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 29175
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 29175
$ws.Range("N88").Value = -29987
$ws.Range("H91").Value = 29175
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 29175
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 29175
$ws.Range("N91").Value = -31983

